$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The document originally starts with the "List of References:" paragraph.
# We need to:
#   1) insert 5 new paragraphs BEFORE it (bold title, interview excerpt,
#      in-text citation, and two blank separator paragraphs), and
#   2) insert 1 new paragraph AFTER it (the BBC Radio 1 reference entry).
#
# Repeatedly calling InsertParagraphBefore() on the same (stale) Range
# object keeps inserting new *empty* paragraphs immediately in front of the
# "List of References:" paragraph, pushing it further down each time - so we
# first create all the empty paragraph slots we need, then go back and fill
# each one in by its now-fixed paragraph index.
# ---------------------------------------------------------------------------
$refsPara = $d.Paragraphs(1)
$refsRange = $refsPara.Range
$refsRange.InsertParagraphBefore()   # slot 1: bold title
$refsRange.InsertParagraphBefore()   # slot 2: interview excerpt
$refsRange.InsertParagraphBefore()   # slot 3: "(BBC Radio 1, 2025)"
$refsRange.InsertParagraphBefore()   # slot 4: blank
$refsRange.InsertParagraphBefore()   # slot 5: blank
# "List of References:" is now paragraph 6.

# ---------------------------------------------------------------------------
# 1) Bold heading paragraph with the video title.
# ---------------------------------------------------------------------------
$titleRange = $d.Paragraphs(1).Range
$titleRange.Font.Bold = $true
$titleRange.Font.BoldBi = $true
$titleRange.InsertBefore("Taylor Swift brings bread & a wedding invite to Greg James | The Life of a Showgirl Interview")

# ---------------------------------------------------------------------------
# 2) Interview excerpt paragraph (Greg / Taylor dialogue separated by a
#    manual line break, [char]11).
# ---------------------------------------------------------------------------
$quoteRange = $d.Paragraphs(2).Range
$quoteText = "Greg: Have you got the itch for another [tour]; something soon?" + [char]11 + "Taylor: No. No, I" + [char]0x2019 + "m just going be really honest with you. Like, I am so tired. "
$quoteRange.InsertBefore($quoteText)

# ---------------------------------------------------------------------------
# 3) Short in-text citation paragraph.
# ---------------------------------------------------------------------------
$citeRange = $d.Paragraphs(3).Range
$citeRange.InsertBefore("(BBC Radio 1, 2025)")

# Paragraphs 4 and 5 stay blank on purpose.

# ---------------------------------------------------------------------------
# 5) New reference-list entry for the BBC Radio 1 video, inserted right
#    after the "List of References:" paragraph (now #6), formatted with a
#    hanging indent and Times New Roman (matching the other reference
#    entries).
# ---------------------------------------------------------------------------
$d.Paragraphs(6).Range.InsertParagraphAfter()
$newRefRange = $d.Paragraphs(7).Range
$newRefRange.ParagraphFormat.LeftIndent = 36
$newRefRange.ParagraphFormat.FirstLineIndent = -36
$newRefRange.Font.Name = "Times New Roman"
$newRefRange.Font.NameFarEast = "Times New Roman"
$newRefRange.Font.NameOther = "Times New Roman"
$newRefRange.Font.NameBi = "Times New Roman"
$newRefRange.Font.BoldBi = $true
$newRefRange.Font.Kerning = 0
$newRefRange.Font.Ligatures = 0
$newRefRange.InsertBefore("BBC Radio 1 (2025, October 3). Taylor Swift brings bread & a wedding invite to Greg James | The Life of a Showgirl Interview [Video]. Youtube. https://www.youtube.com/watch?v=Cqw8fZhggbQ")

# Re-apply the Times New Roman family pieces to just the "Youtube" word so it
# becomes its own run (mirroring the proofed/spell-checked source), while
# keeping the same formatting as its neighbouring runs.
$f1 = $d.Content
$f1.Find.Execute("Youtube")
$f1.Font.Name = "Times New Roman"

$f2 = $d.Content
$f2.Find.Execute("Youtube")
$f2.Font.NameFarEast = "Times New Roman"

$f3 = $d.Content
$f3.Find.Execute("Youtube")
$f3.Font.NameBi = "Times New Roman"

Write-Output "done"
